# The "progression_schema" XML-mapped table now resolves missing/blank
# advancement values to the literal text "none" instead of leaving the
# cell empty (matches the updated XML source/schema under
# D:\Riku\Koulu\Ohjelmointyo\tietokanta\schemas\progression_schema.xml).
# Fill every currently-blank data cell inside the mapped table range
# (B2:K21) with "none".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 21; $r++) {
    for ($c = 2; $c -le 11; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -eq "" -or $v -eq $null) {
            $cell.Value = "none"
        }
    }
}
